$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact / No display for ContactDetail" row (row 11),
# keeping one row at position 10 which becomes "Jurisdiction / United States of America"
$ws.Rows.Item(11).Delete()

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive row (now row 14 after the deletion) gets the literal text "true".
# A plain Value assignment of "true" is auto-coerced to a Boolean by the engine
# (mirrors Excel's smart-typing), so stage it as text in a scratch cell using a
# leading apostrophe (forces text), then paste-special the value through so the
# destination cell keeps the string type without inheriting the quote-prefix style.
$scratch = $ws.Range("A100")
$scratch.Value = "'true"
$scratch.Copy()
$ws.Range("B14").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false
